# Updates the multiplication answer table with newly generated problems.
$d = $word.ActiveDocument

$replacements = @(
    @("87×86=7482", "16×77=1232"),
    @("85×91=7735", "60×49=2940"),
    @("56×35=1960", "71×96=6816"),
    @("41×13=533",  "18×12=216"),
    @("82×17=1394", "86×85=7310"),
    @("35×85=2975", "49×55=2695"),
    @("60×14=840",  "79×44=3476"),
    @("31×66=2046", "94×39=3666"),
    @("66×53=3498", "83×15=1245"),
    @("82×69=5658", "79×96=7584"),
    @("82×43=3526", "19×97=1843"),
    @("23×44=1012", "56×88=4928"),
    @("66×13=858",  "29×88=2552"),
    @("33×20=660",  "73×88=6424"),
    @("48×92=4416", "94×26=2444"),
    @("81×36=2916", "50×69=3450"),
    @("51×74=3774", "58×44=2552"),
    @("82×84=6888", "78×26=2028"),
    @("12×22=264",  "94×50=4700"),
    @("48×53=2544", "41×89=3649"),
    @("61×91=5551", "32×49=1568"),
    @("99×63=6237", "48×74=3552"),
    @("81×98=7938", "18×57=1026"),
    @("85×19=1615", "82×75=6150"),
    @("42×93=3906", "84×97=8148")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
